$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('B36').Value = 'fio8.sh'
$ws.Range('D36').Value = 'd436cb0e1a476937a5a1957ea19a530a'
$ws.Range('E36').Value = 0

$ws.Range('B37').Value = 'fio5.sh'
$ws.Range('D37').Value = '9040171748b12ce1ef6a8536ac982b3e'
$ws.Range('E37').Value = 0

$ws.Range('B38').Value = 'fio2.sh'
$ws.Range('D38').Value = '512b9068b4cd702237ab8242caceab38'
$ws.Range('E38').Value = 0

$ws.Range('B39').Value = 'fio9.sh'
$ws.Range('D39').Value = '61183c5e98db4661b70386d3bfc25ad1'
$ws.Range('E39').Value = 0

$ws.Range('B40').Value = 'fio4.sh'
$ws.Range('D40').Value = 'fcc11badd30b5102a26311a41b951a3c'
$ws.Range('E40').Value = 0

$ws.Range('B41').Value = 'fio7.sh'
$ws.Range('D41').Value = '93d3ec88c2682a6b17d0a0b845c0b772'
$ws.Range('E41').Value = 0

$ws.Range('B43').Value = 'fio3.sh'
$ws.Range('D43').Value = 'f995d8a0bc8d1f2528220691f5156c37'
$ws.Range('E43').Value = 0

$ws.Range('B44').Value = 'fio16.sh'
$ws.Range('D44').Value = '172fa5348233549194189081d49f9f6d'
$ws.Range('E44').Value = 0

$ws.Range('B45').Value = 'fio21.sh'
$ws.Range('D45').Value = '619f1fbf4c45c242ba377bf6323dbc3b'
$ws.Range('E45').Value = 0

$ws.Range('B46').Value = 'fio13.sh'
$ws.Range('D46').Value = '020e043234bba88e36b5941af15254d5'
$ws.Range('E46').Value = 0

$ws.Range('B47').Value = 'fio14.sh'
$ws.Range('D47').Value = '45bac4192c61cfeba14a64386a4b6e65'
$ws.Range('E47').Value = 0

$ws.Range('B48').Value = 'fio12.sh'
$ws.Range('D48').Value = '5dd1be85dd4ab61a7dda604d81f0fc65'
$ws.Range('E48').Value = 0

$ws.Range('B49').Value = 'fio20.sh'
$ws.Range('D49').Value = 'd3184ae79e15874650c264c7fc201803'
$ws.Range('E49').Value = 0

$ws.Range('B50').Value = 'fio18.sh'
$ws.Range('D50').Value = '2709ab002720bf259dfd7f79bf2a2677'
$ws.Range('E50').Value = 0

$ws.Range('B51').Value = 'fio19.sh'
$ws.Range('D51').Value = '15f9d8fcfcd6ec32a51ab99a9e64f19e'
$ws.Range('E51').Value = 0

$ws.Range('B52').Value = 'fio17.sh'
$ws.Range('D52').Value = '5017eef8adf42f9837bbd372d66e547d'
$ws.Range('E52').Value = 0

$ws.Range('B53').Value = 'fio15.sh'
$ws.Range('D53').Value = '0e7fefd4bf12e3b10921488a1733c871'
$ws.Range('E53').Value = 0

$ws.Range('B54').Value = 'fio11.sh'
$ws.Range('D54').Value = 'e22aa548fea442a35b145881113905b5'
$ws.Range('E54').Value = 0

$ws.Range('B55').Value = 'fio21.sh'
$ws.Range('D55').Value = '619f1fbf4c45c242ba377bf6323dbc3b'
$ws.Range('E55').Value = 0

$ws.Range('B56').Value = 'fio24.sh'
$ws.Range('D56').Value = '9889b12792dc5156d8b0bea1b6a9cc63'
$ws.Range('E56').Value = 0

$ws.Range('B57').Value = 'fio12.sh'
$ws.Range('D57').Value = '5dd1be85dd4ab61a7dda604d81f0fc65'
$ws.Range('E57').Value = 0

$ws.Range('B58').Value = 'fio20.sh'
$ws.Range('D58').Value = 'd3184ae79e15874650c264c7fc201803'
$ws.Range('E58').Value = 0

$ws.Range('B59').Value = 'fio23.sh'
$ws.Range('D59').Value = 'de9caef7e5395d10cc8e56ee6c4e6d96'
$ws.Range('E59').Value = 0

$ws.Range('B60').Value = 'run4.txt'
$ws.Range('D60').Value = 'd41d8cd98f00b204e9800998ecf8427e'
$ws.Range('E60').Value = 1

$ws.Range('B61').Value = 'run22.txt'
$ws.Range('D61').Value = 'd41d8cd98f00b204e9800998ecf8427e'
$ws.Range('E61').Value = 1

$ws.Range('B62').Value = 'run23.txt'
$ws.Range('D62').Value = 'd41d8cd98f00b204e9800998ecf8427e'
$ws.Range('E62').Value = 1

$ws.Range('B63').Value = 'run3.txt'
$ws.Range('D63').Value = 'd41d8cd98f00b204e9800998ecf8427e'
$ws.Range('E63').Value = 1

$ws.Range('B64').Value = 'run10.txt'
$ws.Range('D64').Value = 'd41d8cd98f00b204e9800998ecf8427e'
$ws.Range('E64').Value = 1

$ws.Range('B65').Value = 'run8.txt'
$ws.Range('D65').Value = 'd41d8cd98f00b204e9800998ecf8427e'
$ws.Range('E65').Value = 1

$ws.Range('B66').Value = 'run18.txt'
$ws.Range('D66').Value = 'd41d8cd98f00b204e9800998ecf8427e'
$ws.Range('E66').Value = 1

$ws.Range('B67').Value = 'run20.txt'
$ws.Range('D67').Value = 'd41d8cd98f00b204e9800998ecf8427e'
$ws.Range('E67').Value = 1

$ws.Range('B68').Value = 'run2.txt'
$ws.Range('D68').Value = 'd41d8cd98f00b204e9800998ecf8427e'
$ws.Range('E68').Value = 1

$ws.Range('B69').Value = 'run14.txt'
$ws.Range('D69').Value = 'd41d8cd98f00b204e9800998ecf8427e'
$ws.Range('E69').Value = 1

$ws.Range('B70').Value = 'run12.txt'
$ws.Range('D70').Value = 'd41d8cd98f00b204e9800998ecf8427e'
$ws.Range('E70').Value = 1

$ws.Range('B71').Value = 'run13.txt'
$ws.Range('D71').Value = 'd41d8cd98f00b204e9800998ecf8427e'
$ws.Range('E71').Value = 1

$ws.Range('B72').Value = 'run24.txt'
$ws.Range('D72').Value = 'd41d8cd98f00b204e9800998ecf8427e'
$ws.Range('E72').Value = 1

$ws.Range('B73').Value = 'run21.txt'
$ws.Range('D73').Value = 'd41d8cd98f00b204e9800998ecf8427e'
$ws.Range('E73').Value = 1

$ws.Range('B74').Value = 'run11.txt'
$ws.Range('D74').Value = 'd41d8cd98f00b204e9800998ecf8427e'
$ws.Range('E74').Value = 1

$ws.Range('B75').Value = 'run6.txt'
$ws.Range('D75').Value = 'd41d8cd98f00b204e9800998ecf8427e'
$ws.Range('E75').Value = 1

$ws.Range('B76').Value = 'run7.txt'
$ws.Range('D76').Value = 'd41d8cd98f00b204e9800998ecf8427e'
$ws.Range('E76').Value = 1

$ws.Range('B77').Value = 'run17.txt'
$ws.Range('D77').Value = 'd41d8cd98f00b204e9800998ecf8427e'
$ws.Range('E77').Value = 1

$ws.Range('B78').Value = 'run19.txt'
$ws.Range('D78').Value = 'd41d8cd98f00b204e9800998ecf8427e'
$ws.Range('E78').Value = 1

$ws.Range('B79').Value = 'run9.txt'
$ws.Range('D79').Value = 'd41d8cd98f00b204e9800998ecf8427e'
$ws.Range('E79').Value = 1

$ws.Range('B80').Value = 'run15.txt'
$ws.Range('D80').Value = 'd41d8cd98f00b204e9800998ecf8427e'
$ws.Range('E80').Value = 1

$ws.Range('B81').Value = 'run16.txt'
$ws.Range('D81').Value = 'd41d8cd98f00b204e9800998ecf8427e'
$ws.Range('E81').Value = 1

$ws.Range('B82').Value = 'run5.txt'
$ws.Range('D82').Value = 'd41d8cd98f00b204e9800998ecf8427e'
$ws.Range('E82').Value = 1

$ws.Range('B83').Value = 'device.ios2'
$ws.Range('D83').Value = '51b08f41de06d5bd0ff0259c5f0efc3f'
$ws.Range('E83').Value = 0

$ws.Range('B84').Value = 'device.ios3'
$ws.Range('D84').Value = 'bf6d291ed6e2c0330b061c3ce2422e24'
$ws.Range('E84').Value = 0

$ws.Range('B89').Value = 'devlist.vm4_4kall.bkp'
$ws.Range('D89').Value = 'c17a73aa1a4d93f49861b304eedbecae'
$ws.Range('E89').Value = 0

$ws.Range('B91').Value = 'devlist.vm3_4kall.bkp'
$ws.Range('D91').Value = '4fdde708a5376fc76f2391c330a9310c'
$ws.Range('E91').Value = 1

$ws.Range('B96').Value = 'KKD geo tagging (3).xlsx'
$ws.Range('D96').Value = '66c2e9c8ad2fdbda9413f40d1300c29b'
$ws.Range('E96').Value = 1

$ws.Range('B97').Value = 'KKD geo tagging.xlsx'
$ws.Range('D97').Value = '66c2e9c8ad2fdbda9413f40d1300c29b'
$ws.Range('E97').Value = 1

$ws.Range('B98').Value = 'LS____ (2).pdf'
$ws.Range('D98').Value = 'dfa219b9bbb5df3fb8174665ff7c50ae'
$ws.Range('E98').Value = 0

$ws.Range('B99').Value = 'LS____ (1).pdf'
$ws.Range('D99').Value = 'dce5816ff3c5d827355d0edc6e5d9c80'
$ws.Range('E99').Value = 0
